$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.18
$ws.Range("C2").Value = 0.5933333333333334
$ws.Range("J2").Value = 0.003333333333333334
$ws.Range("P2").Value = 0.1666666666666667
$ws.Range("S2").Value = 0.05666666666666666
$ws.Range("B3").Value = 0.01507537688442211
$ws.Range("C3").Value = 0.06030150753768844
$ws.Range("J3").Value = 0.02010050251256281
$ws.Range("P3").Value = 0.7688442211055276
$ws.Range("S3").Value = 0.135678391959799
$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("P4").Value = 0.7413793103448276
$ws.Range("S4").Value = 0.1896551724137931
$ws.Range("B6").Value = 0.07894736842105263
$ws.Range("D6").Value = 0.01754385964912281
$ws.Range("E6").Value = 0.004385964912280702
$ws.Range("F6").Value = 0.05701754385964912
$ws.Range("J6").Value = 0.2280701754385965
$ws.Range("O6").Value = 0.05701754385964912
$ws.Range("Q6").Value = 0.1403508771929824
$ws.Range("R6").Value = 0.1052631578947368
$ws.Range("S6").Value = 0.3114035087719298
$ws.Range("B7").Value = 0.08040201005025126
$ws.Range("D7").Value = 0.02512562814070352
$ws.Range("F7").Value = 0.04020100502512563
$ws.Range("J7").Value = 0.1407035175879397
$ws.Range("O7").Value = 0.02010050251256281
$ws.Range("Q7").Value = 0.1758793969849246
$ws.Range("R7").Value = 0.06532663316582915
$ws.Range("S7").Value = 0.4522613065326633
$ws.Range("B8").Value = 0.07708779443254818
$ws.Range("D8").Value = 0.01927194860813704
$ws.Range("F8").Value = 0.03640256959314775
$ws.Range("J8").Value = 0.1177730192719486
$ws.Range("O8").Value = 0.03426124197002142
$ws.Range("Q8").Value = 0.2077087794432548
$ws.Range("R8").Value = 0.09635974304068523
$ws.Range("S8").Value = 0.411134903640257
$ws.Range("B9").Value = 0.07792207792207792
$ws.Range("D9").Value = 0.02597402597402598
$ws.Range("F9").Value = 0.0735930735930736
$ws.Range("J9").Value = 0.09956709956709957
$ws.Range("O9").Value = 0.01731601731601732
$ws.Range("Q9").Value = 0.1861471861471861
$ws.Range("R9").Value = 0.08658008658008658
$ws.Range("S9").Value = 0.4329004329004329
$ws.Range("B10").Value = 0.1018582243633861
$ws.Range("D10").Value = 0.02408809359944942
$ws.Range("E10").Value = 0.002752924982794219
$ws.Range("F10").Value = 0.0653819683413627
$ws.Range("J10").Value = 0.1183757742601514
$ws.Range("O10").Value = 0.03097040605643496
$ws.Range("Q10").Value = 0.2222986923606332
$ws.Range("R10").Value = 0.08121128699242945
$ws.Range("S10").Value = 0.3530626290433586
$ws.Range("G11").Value = 0.1471571906354515
$ws.Range("J11").Value = 0.1003344481605351
$ws.Range("K11").Value = 0.1839464882943144
$ws.Range("L11").Value = 0.5585284280936454
$ws.Range("S11").Value = 0.01003344481605351
$ws.Range("G12").Value = 0.6982248520710059
$ws.Range("J12").Value = 0.1952662721893491
$ws.Range("K12").Value = 0.02366863905325444
$ws.Range("L12").Value = 0.02366863905325444
$ws.Range("S12").Value = 0.05917159763313609
$ws.Range("G13").Value = 0.7454545454545455
$ws.Range("J13").Value = 0.2545454545454545
$ws.Range("F15").Value = 0.02333333333333333
$ws.Range("H15").Value = 0.1366666666666667
$ws.Range("I15").Value = 0.07000000000000001
$ws.Range("J15").Value = 0.3733333333333334
$ws.Range("K15").Value = 0.06333333333333334
$ws.Range("M15").Value = 0.01333333333333333
$ws.Range("O15").Value = 0.06
$ws.Range("S15").Value = 0.26
$ws.Range("F16").Value = 0.02916666666666667
$ws.Range("H16").Value = 0.1583333333333333
$ws.Range("I16").Value = 0.0625
$ws.Range("J16").Value = 0.4625
$ws.Range("K16").Value = 0.05
$ws.Range("M16").Value = 0.04166666666666666
$ws.Range("O16").Value = 0.08749999999999999
$ws.Range("S16").Value = 0.1083333333333333
$ws.Range("F17").Value = 0.01142857142857143
$ws.Range("H17").Value = 0.1828571428571429
$ws.Range("I17").Value = 0.1085714285714286
$ws.Range("J17").Value = 0.4095238095238095
$ws.Range("K17").Value = 0.07809523809523809
$ws.Range("M17").Value = 0.01523809523809524
$ws.Range("O17").Value = 0.06857142857142857
$ws.Range("S17").Value = 0.1257142857142857
$ws.Range("F18").Value = 0.03181818181818181
$ws.Range("H18").Value = 0.1636363636363636
$ws.Range("I18").Value = 0.1181818181818182
$ws.Range("J18").Value = 0.4318181818181818
$ws.Range("K18").Value = 0.1090909090909091
$ws.Range("M18").Value = 0.004545454545454545
$ws.Range("O18").Value = 0.05454545454545454
$ws.Range("S18").Value = 0.08636363636363636
$ws.Range("F19").Value = 0.02098950524737631
$ws.Range("H19").Value = 0.1971514242878561
$ws.Range("I19").Value = 0.08395802098950525
$ws.Range("J19").Value = 0.3958020989505248
$ws.Range("K19").Value = 0.1026986506746627
$ws.Range("M19").Value = 0.02698650674662669
$ws.Range("N19").Value = 0.0007496251874062968
$ws.Range("O19").Value = 0.07271364317841079
$ws.Range("S19").Value = 0.09895052473763119
